# Append " (Changed main)" to the end of the first paragraph, as three
# separate, distinct runs (" (", "Changed main", ")") following the
# existing "This is a Microsoft word document." run, per the target diff:
#
#   <w:p>
#     <w:r><w:t>This is a Microsoft word document.</w:t></w:r>
#     <w:r><w:t xml:space="preserve"> (</w:t></w:r>
#     <w:r><w:t>Changed main</w:t></w:r>
#     <w:r><w:t>)</w:t></w:r>
#   </w:p>
#
# Plain Range.InsertAfter() calls at the end of the paragraph get
# silently coalesced into the preceding run when the run formatting
# matches (the host mimics Word's own "merge same-format adjacent runs"
# behaviour). To force a genuinely new <w:r> element we instead grow the
# text in a temporary paragraph appended right after paragraph 1, then
# delete the intervening paragraph-mark to splice that text back onto
# paragraph 1 as its own run.

$d = $word.ActiveDocument

function Append-AsNewRun {
    param(
        [int]$ParaIndex,
        [string]$Text
    )

    $para = $d.Paragraphs($ParaIndex)
    $paraRange = $para.Range
    # Exclude the trailing paragraph mark from the range.
    $paraRange.MoveEnd(1, -1) | Out-Null
    $endOfText = $paraRange.End

    # Split the paragraph in two at the end of its text, producing a
    # fresh (empty) paragraph right after it.
    $splitPoint = $d.Range($endOfText, $endOfText)
    $splitPoint.InsertParagraphAfter()

    # Fill the new paragraph with the requested text.
    $newParaIndex = $ParaIndex + 1
    $newPara = $d.Paragraphs($newParaIndex)
    $newParaRange = $newPara.Range
    $newParaRange.MoveEnd(1, -1) | Out-Null
    $newParaRange.InsertAfter($Text)

    # Re-join: delete the paragraph mark that now separates the original
    # paragraph from the new one, merging the new text back in as its
    # own trailing run instead of letting it fuse with the previous run.
    $origPara = $d.Paragraphs($ParaIndex)
    $origEnd = $origPara.Range.End
    $markStart = $origEnd - 1
    $markRange = $d.Range($markStart, $origEnd)
    $markRange.Delete()
}

Append-AsNewRun 1 " ("
Append-AsNewRun 1 "Changed main"
Append-AsNewRun 1 ")"

Write-Output $d.Paragraphs(1).Range.Text
